$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Range("B4").Value = "musterstrasse"
$ws.Range("B8").Select()
